$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.852.94'
$ws.Range("E2").Value = '  -1.72%  '
$ws.Range("D3").Value = '1.833.36'
$ws.Range("E3").Value = '  -1.64%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6937'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07696'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.62%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3051'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.32'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07792'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("B12").Value = 'Litecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '93.30'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.02%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.836.79'
$ws.Range("E13").Value = '  -2.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.096'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.98%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6799'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.57%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.453'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008287'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.15%  '
$ws.Range("D18").Value = '28.882.83'
$ws.Range("E18").Value = '  -1.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.39'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = '2.073.94'
$ws.Range("E20").Value = '  -2.38%  '
$ws.Range("E21").Value = '  -2.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9995'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.449'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1480'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.75%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.64'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.795'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.27'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.539'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.218'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.155'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.187'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05093'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.93%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7730'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.32%  '
$ws.Range("E35").Value = '  -1.76%  '
$ws.Range("E36").Value = '  -3.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.692'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.70%  '
$ws.Range("E38").Value = '  -0.87%  '
$ws.Range("D39").Value = '1.229.83'
$ws.Range("E39").Value = '  -3.73%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.698'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.81%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9515'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '107.98'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.944'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.0000'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.652'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.75%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '1.976.27'
$ws.Range("E46").Value = '  -2.19%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5157'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.30%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '63.98'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -8.76%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.742'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.88%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00000000117'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.926'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.32%  '
